$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 1.29
$ws.Range("N3").Value = 21
$ws.Range("U3").Value = 1.75
$ws.Range("V3").Value = 2
$ws.Range("W3").Value = 11
$ws.Range("Z3").Value = 9
$ws.Range("AD3").Value = 12
$ws.Range("AG3").Value = 151
$ws.Range("AJ3").Value = 23
$ws.Range("AM3").Value = 41
$ws.Range("AN3").Value = 3.6
$ws.Range("AV3").Value = 41
$ws.Range("BC3").Value = 401
$ws.Range("G4").Value = 1.67
$ws.Range("M4").Value = 1.05
$ws.Range("O4").Value = 1.37
$ws.Range("G5").Value = 2.75
$ws.Range("I5").Value = 2.75
$ws.Range("M5").Value = 1.08
$ws.Range("O5").Value = 1.5
$ws.Range("P5").Value = 2.37
$ws.Range("G6").Value = 2.8
$ws.Range("H6").Value = 2.75
$ws.Range("M6").Value = 1.11
$ws.Range("O6").Value = 1.63
$ws.Range("G14").Value = 3.35
$ws.Range("H14").Value = 3.95
$ws.Range("I14").Value = 1.9
$ws.Range("J14").Value = 3.65
$ws.Range("W14").Value = 14.5
$ws.Range("X14").Value = 21
$ws.Range("Y14").Value = 11.75
$ws.Range("AH14").Value = 11
$ws.Range("AL14").Value = 13.5
$ws.Range("AN14").Value = 5.6
$ws.Range("AO14").Value = 17
$ws.Range("AZ14").Value = 30
$ws.Range("O18").Value = 1.22
$ws.Range("P18").Value = 4
$ws.Range("G21").Value = 1.7
$ws.Range("H21").Value = 3.5
$ws.Range("I21").Value = 4.5
$ws.Range("J21").Value = 2.4
$ws.Range("K21").Value = 2.1
$ws.Range("L21").Value = 5
$ws.Range("M21").Value = 1.07
$ws.Range("N21").Value = 9
$ws.Range("O21").Value = 1.33
$ws.Range("P21").Value = 3.25
$ws.Range("Q21").Value = 2.05
$ws.Range("R21").Value = 1.75
$ws.Range("S21").Value = 1.44
$ws.Range("T21").Value = 2.63
$ws.Range("U21").Value = 1.91
$ws.Range("V21").Value = 1.8
$ws.Range("W21").Value = 6.5
$ws.Range("X21").Value = 7.5
$ws.Range("Y21").Value = 8.5
$ws.Range("Z21").Value = 13
$ws.Range("AC21").Value = 9
$ws.Range("AD21").Value = 7
$ws.Range("AE21").Value = 17
$ws.Range("AG21").Value = 351
$ws.Range("AI21").Value = 23
$ws.Range("AJ21").Value = 15
$ws.Range("AK21").Value = 51
$ws.Range("AL21").Value = 41
$ws.Range("AN21").Value = 3.6
$ws.Range("AO21").Value = 9.5
$ws.Range("AQ21").Value = 29
$ws.Range("AT21").Value = 2.63
$ws.Range("AU21").Value = 8.5
$ws.Range("AW21").Value = 6.5
$ws.Range("AX21").Value = 26
$ws.Range("AY21").Value = 34
$ws.Range("AZ21").Value = 101
$ws.Range("BA21").Value = 126
$ws.Range("BB21").Value = 251
$ws.Range("Q28").Value = 2.03
$ws.Range("R28").Value = 1.83
$ws.Range("M33").Value = 1.1
$ws.Range("N33").Value = 7
$ws.Range("O33").Value = 1.5
$ws.Range("P33").Value = 2.5
$ws.Range("Q35").Value = 1.36
$ws.Range("R35").Value = 3.1
$ws.Range("J40").Value = 7.5
$ws.Range("M40").Value = 1.05
$ws.Range("N40").Value = 11
$ws.Range("S40").Value = 1.36
$ws.Range("T40").Value = 3
$ws.Range("U40").Value = 2.1
$ws.Range("V40").Value = 1.67
$ws.Range("Z40").Value = 101
$ws.Range("AS40").Value = 351
$ws.Range("AT40").Value = 3
$ws.Range("AU40").Value = 9.5
$ws.Range("AV40").Value = 67
$ws.Range("AW40").Value = 3.25
$ws.Range("AX40").Value = 6.5
$ws.Range("M41").Value = 1.08
$ws.Range("N41").Value = 8
$ws.Range("Z41").Value = 12
$ws.Range("AB41").Value = 29
$ws.Range("AF41").Value = 51
$ws.Range("AP41").Value = 21
$ws.Range("AS41").Value = 151
$ws.Range("BA41").Value = 126
$ws.Range("BB41").Value = 301
$ws.Range("O42").Value = 1.29
$ws.Range("P42").Value = 3.75
$ws.Range("Q42").Value = 1.9
$ws.Range("R42").Value = 1.95
$ws.Range("M43").Value = 1.06
$ws.Range("N43").Value = 10
$ws.Range("O45").Value = 1.25
$ws.Range("P45").Value = 3.75
$ws.Range("Q45").Value = 1.88
$ws.Range("R45").Value = 1.98
$ws.Range("M46").Value = 1.07
$ws.Range("N46").Value = 9
$ws.Range("Q46").Value = 2.15
$ws.Range("R46").Value = 1.67
$ws.Range("G47").Value = 3.6
$ws.Range("I47").Value = 1.9
$ws.Range("N47").Value = 9.5
$ws.Range("O47").Value = 1.3
$ws.Range("P47").Value = 3.4
$ws.Range("R47").Value = 1.8
$ws.Range("S47").Value = 1.4
$ws.Range("T47").Value = 2.75
$ws.Range("W47").Value = 11
$ws.Range("AA47").Value = 34
$ws.Range("AT47").Value = 2.75
$ws.Range("AZ47").Value = 34
$ws.Range("N55").Value = 9
$ws.Range("Q55").Value = 2.15
$ws.Range("R55").Value = 1.67
$ws.Range("G56").Value = 2.15
$ws.Range("I56").Value = 3.8
$ws.Range("M56").Value = 1.07
$ws.Range("N56").Value = 9
$ws.Range("O56").Value = 1.33
$ws.Range("P56").Value = 3.25
$ws.Range("Q56").Value = 2.1
$ws.Range("R56").Value = 1.7
$ws.Range("X56").Value = 9.5
$ws.Range("Y56").Value = 9
$ws.Range("AL56").Value = 34
$ws.Range("H59").Value = 3.8
$ws.Range("I59").Value = 3.3
$ws.Range("O61").Value = 1.33
$ws.Range("P61").Value = 3.25
$ws.Range("N66").Value = 29
$ws.Range("G68").Value = 1.8
$ws.Range("I68").Value = 4.1
$ws.Range("J68").Value = 2.38
$ws.Range("AQ68").Value = 29
$ws.Range("G83").Value = 1.9
$ws.Range("H83").Value = 3.8
$ws.Range("I83").Value = 3.6
$ws.Range("W83").Value = 9.5
$ws.Range("AO83").Value = 9.5
$ws.Range("Q86").Value = 2.3
$ws.Range("R86").Value = 1.6
$ws.Range("G87").Value = 1.91
$ws.Range("H87").Value = 3.2
$ws.Range("I87").Value = 4.33
$ws.Range("M87").Value = 1.1
$ws.Range("N87").Value = 7
$ws.Range("X87").Value = 7.5
$ws.Range("Z87").Value = 15
$ws.Range("AA87").Value = 19
$ws.Range("AI87").Value = 21
$ws.Range("AK87").Value = 51
$ws.Range("AO87").Value = 11
$ws.Range("AP87").Value = 26
$ws.Range("AR87").Value = 67
$ws.Range("G91").Value = 1.38
$ws.Range("H91").Value = 4.7
$ws.Range("J91").Value = 1.85
$ws.Range("K91").Value = 2.42
$ws.Range("M91").Value = 1.04
$ws.Range("N91").Value = 8.75
$ws.Range("P91").Value = 4
$ws.Range("R91").Value = 2.2
$ws.Range("S91").Value = 1.32
$ws.Range("T91").Value = 3.1
$ws.Range("U91").Value = 1.9
$ws.Range("V91").Value = 1.8
$ws.Range("W91").Value = 7.5
$ws.Range("X91").Value = 6.8
$ws.Range("Z91").Value = 8.75
$ws.Range("AB91").Value = 26
$ws.Range("AC91").Value = 8.75
$ws.Range("AD91").Value = 9.25
$ws.Range("AG91").Value = 700
$ws.Range("AH91").Value = 19.5
$ws.Range("AN91").Value = 3.25
$ws.Range("AO91").Value = 6.2
$ws.Range("AP91").Value = 16
$ws.Range("AQ91").Value = 17
$ws.Range("AR91").Value = 45
$ws.Range("AT91").Value = 3.1
$ws.Range("AU91").Value = 8.5
$ws.Range("AV91").Value = 80
$ws.Range("BA91").Value = 300
